$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A = "Attended the Bootcamp", Column B = "Recommends the Bootcamp",
# Column C = "Result". The Result text was previously computed as
# Attended/Recommends (A/B); it must now be recomputed as
# Recommends/Attended (B/A), formatted like "0.60%" (two decimals with a
# trailing "%" literal appended to the plain ratio, NOT ratio*100).
#
# These are stored as literal text strings (not real percentage numbers),
# so we force text parsing via NumberFormat "@" while writing, then reset
# the cell style back to Normal so the cell doesn't pick up a stray style
# id that wasn't in the source file.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()

    if ($a -eq $null -or $b -eq $null -or $a -eq 0) {
        continue
    }

    $ratio = $b / $a
    $text = "{0:N2}%" -f $ratio

    $cell = $ws.Cells.Item($r, 3)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}
